$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Header: "Website: yogisugandi.site" -> "Website: yogiesugandi.github.io/yogisugandi/"
#    Only the text of the URL run should change - the preceding "Website: "
#    run must stay untouched/separate. A bare Find/Replace across the header
#    tends to coalesce the two adjacent same-format runs into one, so we
#    briefly nudge the URL run's formatting (forcing it to stay a distinct
#    run), perform the text replace, then restore the formatting.
# ---------------------------------------------------------------------------
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(1)
$hdrParas = $hdr.Range.Paragraphs
$hdrLast = $hdrParas.Item($hdrParas.Count)

$siteRng = $hdrLast.Range.Duplicate
if ($siteRng.Find.Execute("yogisugandi.site")) {
    $siteRng.Italic = 1
}

$hdrLast.Range.Find.Execute("yogisugandi.site", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "yogiesugandi.github.io/yogisugandi/", 2)

$newSiteRng = $hdrLast.Range.Duplicate
if ($newSiteRng.Find.Execute("yogiesugandi.github.io/yogisugandi/")) {
    $newSiteRng.Italic = 0
}

# ---------------------------------------------------------------------------
# 2) Body: the "Hormat Saya" (no comma) paragraph - preceded by two empty
#    paragraphs - loses one of those empty paragraphs, gains a first-line
#    indent, and its leading tab character is removed.
# ---------------------------------------------------------------------------
$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("`tHormat Saya")) {
        $targetIdx = $i
    }
}

if ($targetIdx -gt 1) {
    $prev = $d.Paragraphs.Item($targetIdx - 1)
    if ($prev.Range.Text -eq "`r") {
        $prev.Range.Delete()
    }

    # Index shifted by the deletion above - locate the paragraph again.
    $targetIdx = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith("`tHormat Saya")) {
            $targetIdx = $i
        }
    }

    $target = $d.Paragraphs.Item($targetIdx)
    $target.FirstLineIndent = 22.5
    $target.Range.Find.Execute([char]9, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 2)
}
